$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.55'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.91%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '29.95'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '10.44%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.172'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.87%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05718'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.09%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.614'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2.14%'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8600'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '4.67%'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8695'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '3.45%'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1363'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.59%'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07092'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.34%'
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02865'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-5.09%'
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09399'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.03%'
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001517'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.22%'
$ws.Range("B14").Value = 'CoinExToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04142'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.80%'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006019'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.65%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006224'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.31%'
$ws.Range("B17").Value = 'UpBots'
$ws.Range("C17").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.007491'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '3,765.10%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.479'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.03%'
$ws.Range("B19").Value = 'GateToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.055'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.79%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.266'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.83%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '3.51%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.55%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.145'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-11.79%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.41%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.005100'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '14.27%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.001222'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.01%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '23.43%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03774'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.33%'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005719'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '66.41%'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1070'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.45%'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002100'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-17.03%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009333'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '3.99%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005090'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-4.18%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.03%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07509'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002740'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '3.66%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.03%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.03%'
